# Update the cryptos price list (Price / Volume(1h) columns) to the
# latest scraped values, per the GitHub Actions refresh job.
#
# Price cells that look like plain decimal numbers (e.g. "594.92") are
# written with a leading apostrophe so Excel keeps them as literal text
# (matching the source data, which stores prices as strings - some of
# them, like "67.067.80", use a thousands-separator dot and are never
# valid numbers to begin with). The style is reset to "Normal" right
# after so no stray quote-prefix formatting is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.067.80"
$ws.Range("E2").Value = "  +0.78%  "
$ws.Range("D3").Value = "3.501.90"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'594.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").Value = "'173.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.90%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.598"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.66%  "
$ws.Range("E9").Value = "  +3.52%  "
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("E11").Value = "  -1.14%  "
$ws.Range("D12").Value = "4.108.01"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("D14").Value = "'28.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.25%  "
$ws.Range("D15").Value = "67.066.29"
$ws.Range("E15").Value = "  +0.75%  "
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("D17").Value = "3.548.48"
$ws.Range("E17").Value = "  +1.52%  "
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").Value = "'14.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "'394.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("D21").Value = "'8.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("D22").Value = "'73.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "'0.536"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("E25").Value = "  -2.95%  "
$ws.Range("E26").Value = "  -1.65%  "
$ws.Range("D27").Value = "'10.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.40%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "'1.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("E30").Value = "  -1.15%  "
$ws.Range("E31").Value = "  -3.26%  "
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("D33").Value = "'23.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("D34").Value = "'7.39"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.79%  "
$ws.Range("D35").Value = "'1.65"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.90%  "
$ws.Range("D36").Value = "'163.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.59%  "
$ws.Range("D37").Value = "'0.882"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("E39").Value = "  +1.50%  "
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("E41").Value = "  -0.53%  "
$ws.Range("D42").Value = "'27.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.03%  "
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("D44").Value = "2.810.83"
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("E45").Value = "  +2.18%  "
$ws.Range("D46").Value = "'42.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.13%  "
$ws.Range("D47").Value = "'0.0303"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.93%  "
$ws.Range("D48").Value = "'336.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.41%  "
$ws.Range("D49").Value = "'34.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.34%  "
$ws.Range("E50").Value = "  -1.11%  "
$ws.Range("E51").Value = "  -0.09%  "

Write-Output "Updated cryptos list"
